# Deploy schema and regenerate reports
#
# The "Chart Report" sheet feeds both the 3D pie chart and the 3D bar
# chart via the ranges 'Chart Report'!$A$2:$A$25 / $B$2:$B$25. A fresh
# report run changed three of the file-type counts:
#   docx : 142 -> 143
#   html :  28 -> 32
#   rb   :  13 -> 17
#
# Update the source cells (this is what actually drives the charts -
# both embedded charts read straight from this range), then nudge each
# chart's series back onto its own source range so the charts stay
# bound to the refreshed data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart Report")

$ws.Range("B3").Value = 143
$ws.Range("B5").Value = 32
$ws.Range("B6").Value = 17

# Re-assert each chart's series against its (now updated) source range
# so the embedded charts stay bound to the refreshed figures. Use an
# explicit sheet-qualified formula string (rather than a Range object)
# so the series keeps pointing at 'Chart Report'!$B$2:$B$25 instead of
# collapsing to an unqualified reference.
$charts = $ws.ChartObjects()
for ($i = 1; $i -le $charts.Count; $i++) {
    $chart = $charts.Item($i).Chart
    $series = $chart.SeriesCollection(1)
    $valuesFormula = $series.Formula
    if ($valuesFormula -like "*`$B`$2:`$B`$25*") {
        $series.Values = "='Chart Report'!`$B`$2:`$B`$25"
        $series.XValues = "='Chart Report'!`$A`$2:`$A`$25"
    }
}

$wb.RefreshAll()
